$wb = $excel.ActiveWorkbook

# --- Trend_instructions sheet edits ---
$trend = $wb.Worksheets.Item("Trend_instructions")

# B1: "agg_fuel" -> "CEDS_fuel"
$trend.Range("B1").Value = "CEDS_fuel"

# B2: "coal" -> "coal coke"
$trend.Range("B2").Value = "coal coke"

# Update the active selection on this sheet from F3 to B2
$trend.Range("B2").Select()

$wb.Save()
